# Auto-generated Excel COM-interop script
# Applies numeric updates to ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets
# per the authoritative diff of the workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 20000
$ws.Range("J63").Value = 20000
$ws.Range("L63").Value = 20000
$ws.Range("N63").Value = -21248

$ws.Range("H66").Value = 20000
$ws.Range("J66").Value = 20000
$ws.Range("L66").Value = 60000
$ws.Range("N66").Value = -66240

$ws.Range("H74").Value = 4320.8
$ws.Range("I74").Value = 3866.6667
$ws.Range("J74").Value = 5002
$ws.Range("K74").Value = 3866.6667
$ws.Range("L74").Value = 5002
$ws.Range("M74").Value = -2930.6667
$ws.Range("N74").Value = -6874

$ws.Range("H77").Value = 4320.8
$ws.Range("I77").Value = 3866.6667
$ws.Range("J77").Value = 5002
$ws.Range("K77").Value = 19333.3335
$ws.Range("L77").Value = 25010
$ws.Range("M77").Value = -14653.3335
$ws.Range("N77").Value = -34370

$ws.Range("H132").Value = 37906.105
$ws.Range("I132").Value = 39235.965
$ws.Range("K132").Value = 117707.895
$ws.Range("M132").Value = -115177.895

$ws.Range("H135").Value = 5012.6665
$ws.Range("I135").Value = 4016.2856
$ws.Range("J135").Value = 8500
$ws.Range("K135").Value = 36146.5704
$ws.Range("L135").Value = 76500
$ws.Range("M135").Value = -33611.5704
$ws.Range("N135").Value = -81570

$ws.Range("H137").Value = 19231796
$ws.Range("I137").Value = 23810256
$ws.Range("J137").Value = 2268.3
$ws.Range("K137").Value = 71430768
$ws.Range("L137").Value = 6804.900000000001
$ws.Range("M137").Value = -71428218
$ws.Range("N137").Value = -11904.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1478.7727
$ws.Range("I61").Value = 1078.6342
$ws.Range("J61").Value = 6947.3335
$ws.Range("K61").Value = 1078.6342
$ws.Range("L61").Value = 6947.3335
$ws.Range("M61").Value = -866.6342
$ws.Range("N61").Value = -7371.3335

$ws.Range("H122").Value = 5044.125
$ws.Range("I122").Value = 4190.6
$ws.Range("J122").Value = 6466.6665
$ws.Range("K122").Value = 12571.8
$ws.Range("L122").Value = 19399.9995
$ws.Range("M122").Value = -10121.8
$ws.Range("N122").Value = -24299.9995

$ws.Range("H136").Value = 1478.7727
$ws.Range("I136").Value = 1078.6342
$ws.Range("J136").Value = 6947.3335
$ws.Range("K136").Value = 3235.9026
$ws.Range("L136").Value = 20842.0005
$ws.Range("M136").Value = -685.9025999999999
$ws.Range("N136").Value = -25942.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H99").Value = 3201.0908
$ws.Range("I99").Value = 1978.25
$ws.Range("J99").Value = 3899.8572
$ws.Range("K99").Value = 1978.25
$ws.Range("L99").Value = 3899.8572
$ws.Range("M99").Value = -480.25
$ws.Range("N99").Value = -6895.8572

$ws.Range("H107").Value = 310.9
$ws.Range("I107").Value = 310.9
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 310.9
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1609.1
$ws.Range("N107").ClearContents()

$ws.Range("H126").Value = 3201.0908
$ws.Range("I126").Value = 1978.25
$ws.Range("J126").Value = 3899.8572
$ws.Range("K126").Value = 5934.75
$ws.Range("L126").Value = 11699.5716
$ws.Range("M126").Value = -3464.75
$ws.Range("N126").Value = -16639.5716

$ws.Range("H132").Value = 1798.238
$ws.Range("I132").Value = 1230.6471
$ws.Range("J132").Value = 4210.5
$ws.Range("K132").Value = 3691.9413
$ws.Range("L132").Value = 12631.5
$ws.Range("M132").Value = -1161.9413
$ws.Range("N132").Value = -17691.5

$ws.Range("H134").Value = 3928.0952
$ws.Range("I134").Value = 1823.5834
$ws.Range("J134").Value = 6734.1113
$ws.Range("K134").Value = 5470.7502
$ws.Range("L134").Value = 20202.3339
$ws.Range("M134").Value = -2935.7502
$ws.Range("N134").Value = -25272.3339

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1267.5
$ws.Range("I23").Value = 2595.5
$ws.Range("K23").Value = 7786.5
$ws.Range("M23").Value = -7551.5

$ws.Range("H69").Value = 4440
$ws.Range("I69").Value = 890
$ws.Range("J69").Value = 4762.727
$ws.Range("K69").Value = 2670
$ws.Range("L69").Value = 14288.181
$ws.Range("M69").Value = -1859
$ws.Range("N69").Value = -15910.181

$ws.Range("H72").Value = 4440
$ws.Range("I72").Value = 890
$ws.Range("J72").Value = 4762.727
$ws.Range("K72").Value = 8010
$ws.Range("L72").Value = 42864.543
$ws.Range("M72").Value = -3954
$ws.Range("N72").Value = -50976.543

$ws.Range("H82").Value = 3822.2222
$ws.Range("J82").Value = 4237.5
$ws.Range("L82").Value = 12712.5
$ws.Range("N82").Value = -13524.5

$ws.Range("H85").Value = 3822.2222
$ws.Range("J85").Value = 4237.5
$ws.Range("L85").Value = 12712.5
$ws.Range("N85").Value = -15520.5

$ws.Range("H86").Value = 1220
$ws.Range("I86").Value = 366.66666
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 1099.99998
$ws.Range("L86").Value = 7500
$ws.Range("M86").Value = 86.00001999999995
$ws.Range("N86").Value = -9872

$ws.Range("H89").Value = 1220
$ws.Range("I89").Value = 366.66666
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 3299.99994
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = 2628.00006
$ws.Range("N89").Value = -34356

$ws.Range("H141").Value = 4711.1113
$ws.Range("I141").Value = 5816.6665
$ws.Range("K141").Value = 17449.9995
$ws.Range("M141").Value = -12269.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 28825.25
$ws.Range("J93").Value = 28825.25
$ws.Range("L93").Value = 28825.25
$ws.Range("N93").Value = -32569.25

$ws.Range("H132").Value = 3143.5715
$ws.Range("I132").Value = 2570.423
$ws.Range("J132").Value = 10594.5
$ws.Range("K132").Value = 7711.268999999999
$ws.Range("L132").Value = 31783.5
$ws.Range("M132").Value = -5181.268999999999
$ws.Range("N132").Value = -36843.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3083.0476
$ws.Range("I7").Value = 2083.3333
$ws.Range("J7").Value = 3249.6667
$ws.Range("K7").Value = 2083.3333
$ws.Range("L7").Value = 3249.6667
$ws.Range("M7").Value = -1971.3333
$ws.Range("N7").Value = -3473.6667

$ws.Range("H93").Value = 3722
$ws.Range("I93").Value = 4000
$ws.Range("J93").Value = 3444
$ws.Range("K93").Value = 4000
$ws.Range("L93").Value = 3444
$ws.Range("M93").Value = -2752
$ws.Range("N93").Value = -5940

$ws.Range("H126").Value = 3083.0476
$ws.Range("I126").Value = 2083.3333
$ws.Range("J126").Value = 3249.6667
$ws.Range("K126").Value = 6249.999899999999
$ws.Range("L126").Value = 9749.000100000001
$ws.Range("M126").Value = -3779.999899999999
$ws.Range("N126").Value = -14689.0001

$ws.Range("H132").Value = 4077.5833
$ws.Range("I132").Value = 2711.8635
$ws.Range("J132").Value = 6223.7144
$ws.Range("K132").Value = 8135.5905
$ws.Range("L132").Value = 18671.1432
$ws.Range("M132").Value = -5605.5905
$ws.Range("N132").Value = -23731.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 4093.4285
$ws.Range("I49").Value = 2663.5
$ws.Range("J49").Value = 6000
$ws.Range("K49").Value = 2663.5
$ws.Range("L49").Value = 6000
$ws.Range("M49").Value = -2433.5
$ws.Range("N49").Value = -6460

$ws.Range("H126").Value = 51262.4
$ws.Range("I126").Value = 84170.75
$ws.Range("J126").Value = 1899.875
$ws.Range("K126").Value = 252512.25
$ws.Range("L126").Value = 5699.625
$ws.Range("M126").Value = -250042.25

$ws.Range("H136").Value = 1415.5278
$ws.Range("I136").Value = 958.48
$ws.Range("J136").Value = 2454.2727
$ws.Range("K136").Value = 2875.44
$ws.Range("L136").Value = 7362.8181
$ws.Range("M136").Value = -325.4400000000001
$ws.Range("N136").Value = -12462.8181
